# Updates cryptos list data (price + volume change) per upstream refresh.
# Note: Price (column D) cells are plain text in the source data (e.g. some
# contain thousands separators like "3.136.37" which aren't valid numbers).
# For the D values that DO look like plain numbers (e.g. "215.52"), a
# leading apostrophe is used so Excel stores them as text too, matching
# the original cell type instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.507.74'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '3.140.79'
$ws.Range('E3').Value = '  -1.80%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''215.52'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').Value = '''636.62'
$ws.Range('E6').Value = '  +2.83%  '
$ws.Range('E7').Value = '  +1.90%  '
$ws.Range('D8').Value = '''0.772'
$ws.Range('E8').Value = '  +11.42%  '
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('D10').Value = '3.136.63'
$ws.Range('E10').Value = '  -2.06%  '
$ws.Range('D11').Value = '''0.560'
$ws.Range('E11').Value = '  -1.88%  '
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('E13').Value = '  -1.49%  '
$ws.Range('D14').Value = '''5.35'
$ws.Range('E14').Value = '  +2.19%  '
$ws.Range('D15').Value = '89.238.89'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').Value = '3.712.27'
$ws.Range('E16').Value = '  -3.89%  '
$ws.Range('D17').Value = '''32.29'
$ws.Range('E17').Value = '  -2.65%  '
$ws.Range('D18').Value = '3.145.90'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('E19').Value = '  +4.69%  '
$ws.Range('D20').Value = '''0.0000228'
$ws.Range('E20').Value = '  +18.86%  '
$ws.Range('E21').Value = '  -1.77%  '
$ws.Range('D22').Value = '''426.09'
$ws.Range('E22').Value = '  -2.33%  '
$ws.Range('D23').Value = '''8.45'
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('D24').Value = '''4.93'
$ws.Range('E24').Value = '  -3.84%  '
$ws.Range('E25').Value = '  +4.74%  '
$ws.Range('D26').Value = '''82.04'
$ws.Range('E26').Value = '  +8.83%  '
$ws.Range('D27').Value = '''11.52'
$ws.Range('E27').Value = '  -2.58%  '
$ws.Range('D28').Value = '3.298.97'
$ws.Range('E28').Value = '  -4.23%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  -6.99%  '
$ws.Range('D32').Value = '''4.05'
$ws.Range('E32').Value = '  -2.61%  '
$ws.Range('D33').Value = '''8.19'
$ws.Range('E33').Value = '  -3.95%  '
$ws.Range('D34').Value = '''506.98'
$ws.Range('E34').Value = '  -6.03%  '
$ws.Range('D35').Value = '''0.146'
$ws.Range('E35').Value = '  +15.39%  '
$ws.Range('D36').Value = '''7.01'
$ws.Range('E36').Value = '  +1.66%  '
$ws.Range('D37').Value = '''1.29'
$ws.Range('E37').Value = '  +2.70%  '
$ws.Range('D38').Value = '''1.84'
$ws.Range('E38').Value = '  -2.75%  '
$ws.Range('D39').Value = '''22.18'
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('E43').Value = '  -3.18%  '
$ws.Range('D44').Value = '''0.365'
$ws.Range('E44').Value = '  -3.60%  '
$ws.Range('D45').Value = '''146.02'
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('E46').Value = '  +5.28%  '
$ws.Range('D47').Value = '''43.68'
$ws.Range('E47').Value = '  -2.74%  '
$ws.Range('D48').Value = '''0.0666'
$ws.Range('E48').Value = '  +13.44%  '
$ws.Range('D49').Value = '''164.72'
$ws.Range('E49').Value = '  -5.57%  '
$ws.Range('D50').Value = '''0.724'
$ws.Range('E50').Value = '  +1.92%  '
$ws.Range('D51').Value = '''24.23'
$ws.Range('E51').Value = '  -0.59%  '
